$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.925.50'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.748.44'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.44%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -1.15%  '
$ws.Range('E9').Value = '  -2.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.160'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('E11').Value = '  -1.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.44'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -19.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.231.02'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.51'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.585.06'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('E16').Value = '  -1.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.750.78'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('E18').Value = '  +1.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.79'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '354.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('E21').Value = '  -3.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.538'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.02%  '
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0893'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.03%  '
$ws.Range('E29').Value = '  -3.75%  '
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '169.53'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('E32').Value = '  -4.03%  '
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.84'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.43'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('E37').Value = '  -2.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.979'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.36%  '
$ws.Range('E39').Value = '  +7.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.14'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '325.03'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.91'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.73%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.20'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0588'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.26'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0254'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '134.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.623'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.63%  '
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('E51').Value = '  +0.62%  '
